$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously empty answer cell (row 8 - "Werkt er 50%?")
$ws.Range("F8").Value = "als ik het met draatjes verbind wat minder mooi is werkt het"

# Append extra remark to existing answers (E2, E3)
$ws.Range("E2").Value = 'er word soms "gezeeverd" (over ander zaken gepraad) maar ik werk wel veel thuis'
$ws.Range("E3").Value = "soms praat ik mee maar ik werk wel veel thuis"

# Fill in previously empty answer cell (row 5 - visie op de grote vragen)
$ws.Range("F5").Value = "ik weet hoe ik de problemen die ik nu heb ga oplossen want ik heb all veel tetorials gevold"

# Existing answers stay the same text, no change needed for F6/F7/C9/C10 content
# (their underlying shared-string indices shift only because the string table
# was reordered - the visible text is unchanged)

# Update the active selection to F5 (mirrors the saved selection in the file)
$ws.Range("F5").Select()
